# Adherencia Pasarela de pagos.xlsx - "Mi Cuenta" field mapping update
# Adds new rows (new objects) to the ObjetosMiCuenta sheet and tidies up
# a couple of cosmetic workbook/view attributes.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. ObjetosMiCuenta: insert the new rows
# ---------------------------------------------------------------------
$ws6 = $wb.Worksheets.Item("ObjetosMiCuenta")

# --- Block 1: 3 new "list" rows right after row 36 (listNumeroPlan) ---
$ws6.Range("A37:A39").EntireRow.Insert()
$ws6.Range("A36:E36").Copy()
$ws6.Range("A37:E39").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$ws6.Range("A37").Value2 = "IngresaMiCuentaTigo"
$ws6.Range("B37").Value2 = "list"
$ws6.Range("C37").Value2 = "LineaPruebasJuanca"
$ws6.Range("D37").Value2 = "(((//*[@id='lines'])//a[contains(@href,'')]))[11]"
$ws6.Range("E37").Formula = '=CONCATENATE("public By ",B37,C37,"=By.",IF(ISNUMBER(SEARCH("@id=",D37)),"xpath(""","id("""),D37,""");")'

$ws6.Range("A38").Value2 = "IngresaMiCuentaTigo"
$ws6.Range("B38").Value2 = "list"
$ws6.Range("C38").Value2 = "LineaHibridoMariana"
$ws6.Range("D38").Value2 = "(((((//*[@id='lines'])//a[contains(@href,'')])))//*[contains(text(),'Mariana')])[2]"
$ws6.Range("E38").Formula = '=CONCATENATE("public By ",B38,C38,"=By.",IF(ISNUMBER(SEARCH("@id=",D38)),"xpath(""","id("""),D38,""");")'
$ws6.Rows.Item(38).RowHeight = 30

$ws6.Range("A39").Value2 = "IngresaMiCuentaTigo"
$ws6.Range("B39").Value2 = "list"
$ws6.Range("C39").Value2 = "LineaHibridoAna"
$ws6.Range("D39").Value2 = "(((//*[@id='lines'])//a[contains(@href,'')]))[9]"
$ws6.Range("E39").Formula = '=CONCATENATE("public By ",B39,C39,"=By.",IF(ISNUMBER(SEARCH("@id=",D39)),"xpath(""","id("""),D39,""");")'

# --- Block 2: 6 new rows right after row 41 (txtCambiateYa, shifted from
#     the original row 38) and before "CrearTuCuentaTigo" group ---
$ws6.Range("A42:A47").EntireRow.Insert()
$ws6.Range("A41:E41").Copy()
$ws6.Range("A42:E47").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$ws6.Range("A42").Value2 = "IngresaMiCuentaTigo"
$ws6.Range("B42").Value2 = "txt"
$ws6.Range("C42").Value2 = "Usuario"
$ws6.Range("D42").Value2 = "//*[@id='top_menu_aside']/nav/ul/li/a"
$ws6.Range("E42").Formula = '=CONCATENATE("public By ",B42,C42,"=By.",IF(ISNUMBER(SEARCH("@id=",D42)),"xpath(""","id("""),D42,""");")'

$ws6.Range("A43").Value2 = "IngresaMiCuentaTigo"
$ws6.Range("B43").Value2 = "txt"
$ws6.Range("C43").Value2 = "CerrarSesion"
$ws6.Range("D43").Value2 = "//*[@id='top_menu_aside']/nav/ul/li/ul/li[2]/a"
$ws6.Range("E43").Formula = '=CONCATENATE("public By ",B43,C43,"=By.",IF(ISNUMBER(SEARCH("@id=",D43)),"xpath(""","id("""),D43,""");")'

$ws6.Range("A44").Value2 = "IngresaMiCuentaTigo"
$ws6.Range("B44").Value2 = "btn"
$ws6.Range("C44").Value2 = "CambiarCuenta"
$ws6.Range("D44").Value2 = "//*[@id='addNew']"
$ws6.Range("E44").Formula = '=CONCATENATE("public By ",B44,C44,"=By.",IF(ISNUMBER(SEARCH("@id=",D44)),"xpath(""","id("""),D44,""");")'

$ws6.Range("A45").Value2 = "IngresaMiCuentaTigo"
$ws6.Range("B45").Value2 = "lb"
$ws6.Range("C45").Value2 = "Planes"
$ws6.Range("D45").Value2 = "//*[@id='main-content']/div[2]"
$ws6.Range("E45").Formula = '=CONCATENATE("public By ",B45,C45,"=By.",IF(ISNUMBER(SEARCH("@id=",D45)),"xpath(""","id("""),D45,""");")'

$ws6.Range("A46").Value2 = "IngresaMiCuentaTigo"
$ws6.Range("B46").Value2 = "lb"
$ws6.Range("C46").Value2 = "MejorarPlan"
$ws6.Range("D46").Value2 = "//*[@id='compras-noplan-container']/h3"
$ws6.Range("E46").Formula = '=CONCATENATE("public By ",B46,C46,"=By.",IF(ISNUMBER(SEARCH("@id=",D46)),"xpath(""","id("""),D46,""");")'

$ws6.Range("A47").Value2 = "IngresaMiCuentaTigo"
$ws6.Range("B47").Value2 = "lb"
$ws6.Range("C47").Value2 = "DetallePlan"
$ws6.Range("D47").Value2 = "//*[@id='compras-noplan-container']/p"
$ws6.Range("E47").Formula = '=CONCATENATE("public By ",B47,C47,"=By.",IF(ISNUMBER(SEARCH("@id=",D47)),"xpath(""","id("""),D47,""");")'

# --- Column E is now wider (no longer "best fit") ---
$ws6.Columns.Item(5).ColumnWidth = 113.41666666666666

# --- View: ObjetosMiCuenta becomes the active sheet / tab ---
$ws6.Activate()
$ws6.Range("A52").Select()

$wb.Save()
